$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Cart_Page: fill in the previously-blank rows 7-10 with the new
# "no promotion" / "checkout button" locators, and drop old row 11.
# ---------------------------------------------------------------
$cart = $wb.Worksheets.Item("Cart_Page")

# Row 11 becomes fully empty (no cell at all) in the target file.
$cart.Rows.Item(11).Clear()

$cart.Range("C7").Value = "swal2-modal"
$cart.Range("A7").Value = "no promotion"
$cart.Range("B7").Value = "class name"

$cart.Range("A8").Value = "no promotion ok button"
$cart.Range("B8").Value = "class name"
$cart.Range("C8").Value = "swal2-confirm"

$cart.Range("C9").Value = "checkoutTop-CartPage"
$cart.Range("C10").Value = "checkoutBottom-CartPage"

$cart.Range("A9").Value = "checkout button top"
$cart.Range("A10").Value = "checkout button bottom"

$cart.Range("B9").Value = "id"
$cart.Range("B10").Value = "id"

# ---------------------------------------------------------------
# My_Account_Page: add "add shipping address" locators.
# ---------------------------------------------------------------
$acct = $wb.Worksheets.Item("My_Account_Page")

# Insert two fresh rows at 13 (pushes old 13/14/15 down to 15/16/17).
$acct.Rows.Item(13).Insert()
$acct.Rows.Item(13).Insert()

# New row 14's D/E blanks should carry style 23 (same as row 15's),
# not the style 15 inherited from the insert-above default.
$acct.Range("B15:C15").Copy()
$acct.Range("D14:E14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The inserted row 14 shouldn't carry an I-column cell.
$acct.Range("I14").Clear()

# New cell text is introduced in the same order the original commit
# used, so the shared-string table lands in the same sequence.
$acct.Range("C14").Value = "Shipping_link__300no"
$acct.Range("A14").Value = "add change shipping buttons"

# New trailing rows describing the "Add New Shipping Address" modal.
$acct.Range("A18").Value = "modal"
$acct.Range("C18").Value = "modal"
$acct.Range("A18:C18").VerticalAlignment = -4160
$acct.Range("D18").Value = "Add New Shipping Address"

$acct.Range("A19").Value = "ad"

$acct.Range("A13").Value = "add shipping address button"
$acct.Range("C13").Value = '//*[@id="modal"]/div[1]/div/form/button'

# Drop the now-unused D/E cells on the (shifted) order-table row.
$acct.Range("D17:E17").Clear()

$acct.Range("B14").Value = "class name"
$acct.Range("B13").Value = "xpath"
$acct.Range("B18").Value = "class name"

# ---------------------------------------------------------------
# Sheet/view bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------

# "Page" sheet: selection moved from A7 to A9:C9.
$page = $wb.Worksheets.Item("Page")
$page.Range("A9:C9").Select()

# My_Account_Page becomes the active tab, with a fresh selection.
$acct.Activate()
$acct.Range("A13").Select()
